$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B24 value from 0 to 1 (this will recalculate dependent formulas
# in C24, B1, C1 automatically)
$ws.Range("B24").Value = 1

# Move the selection / active cell to B25 (GoTo without "Altered")
$ws.Range("B25").Select()
